# This workbook is a daily price log for "Albahaca" (Femacal de La Calera).
# The commit turns it into a weekly log by inserting one new weekly record
# right after the header block of existing rows, pushing the previously
# existing row 26 (and everything below it) down by one row, and then
# filling the freshly inserted row 26 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; this shifts old rows 26-92 down to
# rows 27-93 (carrying their values/formatting with them), matching the
# diff where every row from 27 to 93 now holds what used to be in the row
# directly above it, and a brand-new row 93 appears with the data that used
# to be the last row (92).
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44536
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 100112052
$ws.Range("G26").Value = "Albahaca"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 90
$ws.Range("K26").Value = 4000
$ws.Range("L26").Value = 4500
$ws.Range("M26").Value = 4222
$ws.Range("N26").Value = "`$/docena de matas"
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 704
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = "Hortaliza"
